$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.358.57"
$ws.Range("E2").Value = "  -1.20%  "

$ws.Range("D3").Value = "2.050.63"
$ws.Range("E3").Value = "  -1.24%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.52"
$ws.Range("E5").Value = "  -1.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.621"
$ws.Range("E6").Value = "  -0.62%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.10"
$ws.Range("E8").Value = "  -3.25%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.385"
$ws.Range("E9").Value = "  -2.36%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0771"
$ws.Range("E10").Value = "  -2.23%  "

$ws.Range("E11").Value = "  +1.56%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.76"
$ws.Range("E12").Value = "  -0.17%  "

$ws.Range("D13").Value = "2.349.15"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.59"
$ws.Range("E14").Value = "  -2.89%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.756"
$ws.Range("E15").Value = "  -2.37%  "

$ws.Range("E16").Value = "  -1.90%  "

$ws.Range("D17").Value = "2.051.90"
$ws.Range("E17").Value = "  -2.13%  "

$ws.Range("D18").Value = "37.328.39"
$ws.Range("E18").Value = "  -1.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.98"
$ws.Range("E19").Value = "  -2.66%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.77"
$ws.Range("E20").Value = "  -2.47%  "

$ws.Range("D21").Value = "0.0₃0826"
$ws.Range("E21").Value = "  -2.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.43"
$ws.Range("E22").Value = "  -0.79%  "

$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.37"
$ws.Range("E24").Value = "  +0.60%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.33"
$ws.Range("E25").Value = "  -3.59%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.50"
$ws.Range("E26").Value = "  +3.62%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.77"
$ws.Range("E27").Value = "  -1.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.131"
$ws.Range("E28").Value = "  -3.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.19"
$ws.Range("E29").Value = "  -1.60%  "

$ws.Range("E30").Value = "  -4.45%  "

$ws.Range("E31").Value = "  +0.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.54"
$ws.Range("E32").Value = "  -3.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0625"
$ws.Range("E33").Value = "  -1.26%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.58"
$ws.Range("E34").Value = "  -4.03%  "

$ws.Range("E35").Value = "  -0.58%  "

$ws.Range("E36").Value = "  -0.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.29"
$ws.Range("E37").Value = "  -3.45%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.28"
$ws.Range("E39").Value = "  -1.99%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0227"
$ws.Range("E40").Value = "  +5.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.09"
$ws.Range("E41").Value = "  -1.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0954"
$ws.Range("E42").Value = "  -2.20%  "

$ws.Range("E43").Value = "  +0.60%  "

$ws.Range("D44").Value = "1.481.95"
$ws.Range("E44").Value = "  +2.82%  "

$ws.Range("E45").Value = "  +3.26%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.65"
$ws.Range("E46").Value = "  +0.39%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.03"
$ws.Range("E47").Value = "  -3.39%  "

$ws.Range("E48").Value = "  -3.10%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.22"
$ws.Range("E49").Value = "  -2.49%  "

$ws.Range("E50").Value = "  -2.58%  "

$ws.Range("D51").Value = "2.235.91"
$ws.Range("E51").Value = "  -1.42%  "
